$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New diaper log rows (216-241): A=day index, B=date serial, C=total_diaper,
# D=total_wet, E=total_dirty, F=total_both ---

$colA = New-Object 'object[,]' 26,1
$colB = New-Object 'object[,]' 26,1
$colsCDEF = New-Object 'object[,]' 26,4

$colA[0,0] = 215
$colB[0,0] = 44087
$colsCDEF[0,0] = 5
$colsCDEF[0,1] = 4
$colsCDEF[0,2] = 0
$colsCDEF[0,3] = 1
$colA[1,0] = 216
$colB[1,0] = 44088
$colsCDEF[1,0] = 5
$colsCDEF[1,1] = 4
$colsCDEF[1,2] = 0
$colsCDEF[1,3] = 1
$colA[2,0] = 217
$colB[2,0] = 44089
$colsCDEF[2,0] = 5
$colsCDEF[2,1] = 5
$colsCDEF[2,2] = 0
$colsCDEF[2,3] = 0
$colA[3,0] = 218
$colB[3,0] = 44090
$colsCDEF[3,0] = 5
$colsCDEF[3,1] = 5
$colsCDEF[3,2] = 0
$colsCDEF[3,3] = 0
$colA[4,0] = 219
$colB[4,0] = 44091
$colsCDEF[4,0] = 6
$colsCDEF[4,1] = 4
$colsCDEF[4,2] = 0
$colsCDEF[4,3] = 1
$colA[5,0] = 220
$colB[5,0] = 44092
$colsCDEF[5,0] = 5
$colsCDEF[5,1] = 5
$colsCDEF[5,2] = 0
$colsCDEF[5,3] = 0
$colA[6,0] = 221
$colB[6,0] = 44093
$colsCDEF[6,0] = 5
$colsCDEF[6,1] = 5
$colsCDEF[6,2] = 0
$colsCDEF[6,3] = 0
$colA[7,0] = 222
$colB[7,0] = 44094
$colsCDEF[7,0] = 5
$colsCDEF[7,1] = 4
$colsCDEF[7,2] = 0
$colsCDEF[7,3] = 1
$colA[8,0] = 223
$colB[8,0] = 44095
$colsCDEF[8,0] = 5
$colsCDEF[8,1] = 5
$colsCDEF[8,2] = 0
$colsCDEF[8,3] = 0
$colA[9,0] = 224
$colB[9,0] = 44096
$colsCDEF[9,0] = 6
$colsCDEF[9,1] = 6
$colsCDEF[9,2] = 0
$colsCDEF[9,3] = 0
$colA[10,0] = 225
$colB[10,0] = 44097
$colsCDEF[10,0] = 6
$colsCDEF[10,1] = 6
$colsCDEF[10,2] = 0
$colsCDEF[10,3] = 0
$colA[11,0] = 226
$colB[11,0] = 44098
$colsCDEF[11,0] = 6
$colsCDEF[11,1] = 6
$colsCDEF[11,2] = 0
$colsCDEF[11,3] = 0
$colA[12,0] = 227
$colB[12,0] = 44099
$colsCDEF[12,0] = 5
$colsCDEF[12,1] = 5
$colsCDEF[12,2] = 0
$colsCDEF[12,3] = 0
$colA[13,0] = 228
$colB[13,0] = 44100
$colsCDEF[13,0] = 5
$colsCDEF[13,1] = 4
$colsCDEF[13,2] = 0
$colsCDEF[13,3] = 1
$colA[14,0] = 229
$colB[14,0] = 44101
$colsCDEF[14,0] = 6
$colsCDEF[14,1] = 5
$colsCDEF[14,2] = 0
$colsCDEF[14,3] = 1
$colA[15,0] = 230
$colB[15,0] = 44102
$colsCDEF[15,0] = 5
$colsCDEF[15,1] = 5
$colsCDEF[15,2] = 0
$colsCDEF[15,3] = 0
$colA[16,0] = 231
$colB[16,0] = 44103
$colsCDEF[16,0] = 5
$colsCDEF[16,1] = 5
$colsCDEF[16,2] = 0
$colsCDEF[16,3] = 0
$colA[17,0] = 232
$colB[17,0] = 44104
$colsCDEF[17,0] = 6
$colsCDEF[17,1] = 5
$colsCDEF[17,2] = 0
$colsCDEF[17,3] = 1
$colA[18,0] = 233
$colB[18,0] = 44105
$colsCDEF[18,0] = 5
$colsCDEF[18,1] = 5
$colsCDEF[18,2] = 0
$colsCDEF[18,3] = 0
$colA[19,0] = 234
$colB[19,0] = 44106
$colsCDEF[19,0] = 5
$colsCDEF[19,1] = 5
$colsCDEF[19,2] = 0
$colsCDEF[19,3] = 0
$colA[20,0] = 235
$colB[20,0] = 44107
$colsCDEF[20,0] = 5
$colsCDEF[20,1] = 4
$colsCDEF[20,2] = 0
$colsCDEF[20,3] = 1
$colA[21,0] = 236
$colB[21,0] = 44108
$colsCDEF[21,0] = 6
$colsCDEF[21,1] = 5
$colsCDEF[21,2] = 0
$colsCDEF[21,3] = 1
$colA[22,0] = 237
$colB[22,0] = 44109
$colsCDEF[22,0] = 4
$colsCDEF[22,1] = 4
$colsCDEF[22,2] = 0
$colsCDEF[22,3] = 0
$colA[23,0] = 238
$colB[23,0] = 44110
$colsCDEF[23,0] = 7
$colsCDEF[23,1] = 7
$colsCDEF[23,2] = 0
$colsCDEF[23,3] = 0
$colA[24,0] = 239
$colB[24,0] = 44111
$colsCDEF[24,0] = 6
$colsCDEF[24,1] = 5
$colsCDEF[24,2] = 0
$colsCDEF[24,3] = 1
$colA[25,0] = 240
$colB[25,0] = 44112
$colsCDEF[25,0] = 5
$colsCDEF[25,1] = 5
$colsCDEF[25,2] = 0
$colsCDEF[25,3] = 0

$ws.Range("A216:A241").Value = $colA
$ws.Range("B216:B241").Value = $colB
$ws.Range("C216:F241").Value = $colsCDEF

# Apply the same date number format as the existing date column (copy format only)
$ws.Range("B215").Copy()
$ws.Range("B216:B241").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the sheet selection / active cell to reflect the appended data,
# matching the author's final view state.
[void]$ws.Range("A242").Select()
